$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 386, pushing existing rows 386-420 down to 387-421.
$ws.Rows.Item(386).Insert()

# Populate the newly inserted row 386 with a new weekly price record
# (same market/category as the row that used to be at 386, but a new
# date and volume).
$ws.Range("A386").Value = 10
$ws.Range("B386").Value = "Vega Modelo de Temuco"
$ws.Range("C386").Value = "La Araucanía"
$ws.Range("D386").Value = 45166
$ws.Range("E386").Value = 9
$ws.Range("F386").Value = 100112039
$ws.Range("G386").Value = "Ciboulette"
$ws.Range("H386").Value = "Sin especificar"
$ws.Range("I386").Value = "Primera"
$ws.Range("J386").Value = 30
$ws.Range("K386").Value = 7000
$ws.Range("L386").Value = 7000
$ws.Range("M386").Value = 7000
$ws.Range("N386").Value = "$/docena de atados"
$ws.Range("O386").Value = "Provincia de Cautín"
$ws.Range("P386").Value = 2333
$ws.Range("Q386").Value = 3
$ws.Range("R386").Value = "Hortaliza"
